# Generate Report for Handoff
# The two tracked source files (1ea166bf-... and da2bbf5e-...) swap which
# table row they occupy, and their localization status moves forward:
#   da2bbf5e -> now reported in row 2, status "Ready for handoff"
#   1ea166bf -> now reported in row 3, status "Ready for handoff"
#     (with a stale-handback error noted in the Error Detail column)

$wb = $excel.ActiveWorkbook

$file1 = "1ea166bf-e461-453a-ab2c-b97c0ab01afb"
$file2 = "da2bbf5e-40e8-42ee-ae86-290e53dc2164"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01bc330e94a9c1fc19cfbb3f1a9f1633812492ef/e2e/$file1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a02184e68c0b80432011bfe6f009411a1ea6493/e2e/$file1.md."

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$file2.md"
$wsOverview.Range("B2").Value = "e2e\$file2.md"

$wsOverview.Range("A3").Value = "$file1.md"
$wsOverview.Range("B3").Value = "e2e\$file1.md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 02:47:47"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$file2.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$file1.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$file2.md"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("G2").Value = "$file2.7cdc7e7add6c87705a057a3a11ea6eb79bb7ea2f.zh-cn.xlf"
$wsZhCn.Range("I2").Value = "$file2.md"
$wsZhCn.Range("J2").Value = "$file2.7cdc7e7add6c87705a057a3a11ea6eb79bb7ea2f.zh-cn.xlf"

$wsZhCn.Range("A3").Value = "$file1.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "$file1.548ed771fd3587dbc6956a5d9d227c19c3ec2cd2.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 02:47:42"
$wsZhCn.Range("I3").Value = "$file1.md"
$wsZhCn.Range("J3").Value = "$file1.548ed771fd3587dbc6956a5d9d227c19c3ec2cd2.zh-cn.xlf"
$wsZhCn.Range("P3").Value = $errorDetail

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2' -or $addr -eq '$I$2') {
        $hl.TextToDisplay = "$file2.md"
    } elseif ($addr -eq '$A$3' -or $addr -eq '$I$3') {
        $hl.TextToDisplay = "$file1.md"
    }
}

$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$file2.md"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("G2").Value = "$file2.7cdc7e7add6c87705a057a3a11ea6eb79bb7ea2f.de-de.xlf"
$wsDeDe.Range("I2").Value = "$file2.md"
$wsDeDe.Range("J2").Value = "$file2.7cdc7e7add6c87705a057a3a11ea6eb79bb7ea2f.de-de.xlf"

$wsDeDe.Range("A3").Value = "$file1.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "$file1.548ed771fd3587dbc6956a5d9d227c19c3ec2cd2.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 02:47:47"
$wsDeDe.Range("I3").Value = "$file1.md"
$wsDeDe.Range("J3").Value = "$file1.548ed771fd3587dbc6956a5d9d227c19c3ec2cd2.de-de.xlf"
$wsDeDe.Range("P3").Value = $errorDetail

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2' -or $addr -eq '$I$2') {
        $hl.TextToDisplay = "$file2.md"
    } elseif ($addr -eq '$A$3' -or $addr -eq '$I$3') {
        $hl.TextToDisplay = "$file1.md"
    }
}

$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
